$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, pushing existing rows 40-73 down to 41-74
$ws.Rows("40:40").Insert()

# Copy the style of the date cell from the row below (now row 41, formerly row 40) into new D40
$ws.Range("D41").Copy()
$ws.Range("D40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row 40 with data
$ws.Range("A40").Value = 4
$ws.Range("B40").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C40").Value = 'Los Lagos'
$ws.Range("D40").Value = 44484
$ws.Range("E40").Value = 10
$ws.Range("F40").Value = 100112022
$ws.Range("G40").Value = 'Arveja Verde'
$ws.Range("H40").Value = 'Sin especificar'
$ws.Range("I40").Value = 'Primera'
$ws.Range("J40").Value = 90
$ws.Range("K40").Value = 25000
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = 25000
$ws.Range("N40").Value = '$/saco 25 kilos'
$ws.Range("O40").Value = 'Región Metropolitana'
$ws.Range("P40").Value = 1000
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = 'Hortaliza'
